# Update the cryptocurrency price (D) and 1h-volume-change (E) columns
# with freshly scraped values, per the scheduled GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$d2Style = $ws.Range("D2").Style
$ws.Range("D2").Value = "'28.111.92"
$ws.Range("D2").Style = $d2Style
$ws.Range("E2").Value = "  +1.70%  "
$d3Style = $ws.Range("D3").Style
$ws.Range("D3").Value = "'1.791.48"
$ws.Range("D3").Style = $d3Style
$ws.Range("E3").Value = "  +1.91%  "
$ws.Range("E4").Value = "  -0.58%  "
$d5Style = $ws.Range("D5").Style
$ws.Range("D5").Value = "'324.20"
$ws.Range("D5").Style = $d5Style
$ws.Range("E5").Value = "  -0.57%  "
$ws.Range("E6").Value = "  -0.31%  "
$ws.Range("E7").Value = "  -3.23%  "
$d8Style = $ws.Range("D8").Style
$ws.Range("D8").Value = "'0.3633"
$ws.Range("D8").Style = $d8Style
$ws.Range("E8").Value = "  -2.50%  "
$d9Style = $ws.Range("D9").Style
$ws.Range("D9").Value = "'44.62"
$ws.Range("D9").Style = $d9Style
$ws.Range("E9").Value = "  -2.19%  "
$d10Style = $ws.Range("D10").Style
$ws.Range("D10").Value = "'0.07532"
$ws.Range("D10").Style = $d10Style
$ws.Range("E10").Value = "  -3.15%  "
$ws.Range("E11").Value = "  -1.09%  "
$d12Style = $ws.Range("D12").Style
$ws.Range("D12").Value = "'1.000"
$ws.Range("D12").Style = $d12Style
$ws.Range("E12").Value = "  -0.29%  "
$ws.Range("E13").Value = "  -0.02%  "
$d14Style = $ws.Range("D14").Style
$ws.Range("D14").Value = "'6.171"
$ws.Range("D14").Style = $d14Style
$ws.Range("E14").Value = "  -0.38%  "
$d15Style = $ws.Range("D15").Style
$ws.Range("D15").Value = "'7.354"
$ws.Range("D15").Style = $d15Style
$ws.Range("E15").Value = "  -0.25%  "
$d16Style = $ws.Range("D16").Style
$ws.Range("D16").Value = "'1.781.65"
$ws.Range("D16").Style = $d16Style
$ws.Range("E16").Value = "  +1.26%  "
$ws.Range("E17").Value = "  +1.03%  "
$ws.Range("E18").Value = "  -1.23%  "
$d19Style = $ws.Range("D19").Style
$ws.Range("D19").Value = "'0.06345"
$ws.Range("D19").Style = $d19Style
$ws.Range("E19").Value = "  +1.47%  "
$d20Style = $ws.Range("D20").Style
$ws.Range("D20").Value = "'0.9998"
$ws.Range("D20").Style = $d20Style
$ws.Range("E20").Value = "  -0.29%  "
$ws.Range("E21").Value = "  -0.69%  "
$d22Style = $ws.Range("D22").Style
$ws.Range("D22").Value = "'5.965"
$ws.Range("D22").Style = $d22Style
$ws.Range("E22").Value = "  -3.61%  "
$d23Style = $ws.Range("D23").Style
$ws.Range("D23").Value = "'28.117.41"
$ws.Range("D23").Style = $d23Style
$ws.Range("E23").Value = "  +1.53%  "
$ws.Range("E24").Value = "  -1.80%  "
$d25Style = $ws.Range("D25").Style
$ws.Range("D25").Value = "'2.157"
$ws.Range("D25").Style = $d25Style
$ws.Range("E25").Value = "  -7.50%  "
$d26Style = $ws.Range("D26").Style
$ws.Range("D26").Value = "'160.36"
$ws.Range("D26").Style = $d26Style
$ws.Range("E26").Value = "  +4.35%  "
$d27Style = $ws.Range("D27").Style
$ws.Range("D27").Value = "'20.38"
$ws.Range("D27").Style = $d27Style
$ws.Range("E27").Value = "  -2.10%  "
$d28Style = $ws.Range("D28").Style
$ws.Range("D28").Value = "'1.983.71"
$ws.Range("D28").Style = $d28Style
$ws.Range("E28").Value = "  +1.27%  "
$ws.Range("E29").Value = "  -6.95%  "
$d30Style = $ws.Range("D30").Style
$ws.Range("D30").Value = "'126.99"
$ws.Range("D30").Style = $d30Style
$ws.Range("E30").Value = "  -1.56%  "
$ws.Range("E31").Value = "  -3.43%  "
$d32Style = $ws.Range("D32").Style
$ws.Range("D32").Value = "'5.732"
$ws.Range("D32").Style = $d32Style
$ws.Range("E32").Value = "  -0.70%  "
$d33Style = $ws.Range("D33").Style
$ws.Range("D33").Value = "'0.09018"
$ws.Range("D33").Style = $d33Style
$ws.Range("E33").Value = "  -2.52%  "
$d34Style = $ws.Range("D34").Style
$ws.Range("D34").Value = "'3.502"
$ws.Range("D34").Style = $d34Style
$ws.Range("E34").Value = "  -5.14%  "
$d35Style = $ws.Range("D35").Style
$ws.Range("D35").Value = "'12.70"
$ws.Range("D35").Style = $d35Style
$ws.Range("E35").Value = "  -0.54%  "
$ws.Range("E36").Value = "  -0.50%  "
$d37Style = $ws.Range("D37").Style
$ws.Range("D37").Value = "'5.106"
$ws.Range("D37").Style = $d37Style
$ws.Range("E37").Value = "  +0.32%  "
$d38Style = $ws.Range("D38").Style
$ws.Range("D38").Value = "'0.6475"
$ws.Range("D38").Style = $d38Style
$ws.Range("E38").Value = "  -0.27%  "
$ws.Range("E39").Value = "  -3.03%  "
$d40Style = $ws.Range("D40").Style
$ws.Range("D40").Value = "'0.06071"
$ws.Range("D40").Style = $d40Style
$ws.Range("E40").Value = "  -0.85%  "
$ws.Range("E42").Value = "  -0.03%  "
$d43Style = $ws.Range("D43").Style
$ws.Range("D43").Value = "'0.9995"
$ws.Range("D43").Style = $d43Style
$ws.Range("E43").Value = "  -0.30%  "
$d44Style = $ws.Range("D44").Style
$ws.Range("D44").Value = "'7.895"
$ws.Range("D44").Style = $d44Style
$ws.Range("E44").Value = "  -1.52%  "
$d45Style = $ws.Range("D45").Style
$ws.Range("D45").Value = "'13.71"
$ws.Range("D45").Style = $d45Style
$ws.Range("E45").Value = "  -1.23%  "
$d46Style = $ws.Range("D46").Style
$ws.Range("D46").Value = "'0.6002"
$ws.Range("D46").Style = $d46Style
$ws.Range("E46").Value = "  +0.00%  "
$d47Style = $ws.Range("D47").Style
$ws.Range("D47").Value = "'3.711"
$ws.Range("D47").Style = $d47Style
$ws.Range("E47").Value = "  -1.00%  "
$d48Style = $ws.Range("D48").Style
$ws.Range("D48").Value = "'124.57"
$ws.Range("D48").Style = $d48Style
$ws.Range("E48").Value = "  -0.97%  "
$ws.Range("E49").Value = "  -0.27%  "
$d50Style = $ws.Range("D50").Style
$ws.Range("D50").Value = "'1.155"
$ws.Range("D50").Style = $d50Style
$ws.Range("E50").Value = "  +0.84%  "
$d51Style = $ws.Range("D51").Style
$ws.Range("D51").Value = "'0.06956"
$ws.Range("D51").Style = $d51Style
